$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 139, pushing the existing rows 139-141 down to 140-142.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new record's values.
$ws.Cells.Item(139, 1).Value = 4
$ws.Cells.Item(139, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(139, 3).Value = "Los Lagos"
$ws.Cells.Item(139, 4).Value = 44448
$ws.Cells.Item(139, 5).Value = 10
$ws.Cells.Item(139, 6).Value = 100112043
$ws.Cells.Item(139, 7).Value = "Pepino ensalada"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 200
$ws.Cells.Item(139, 11).Value = 21000
$ws.Cells.Item(139, 12).Value = 21000
$ws.Cells.Item(139, 13).Value = 21000
$ws.Cells.Item(139, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(139, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(139, 16).Value = 350
$ws.Cells.Item(139, 17).Value = 60
$ws.Cells.Item(139, 18).Value = "Hortaliza"
